$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New experiment section: "Experiment #58, Async-Cent" -----------------
$ws.Range("A89").Value = "Experiment #58, Async-Cent"

$ws.Range("A91").Value = "Exchange#"
$ws.Range("B91").Value = "Replica#"
$ws.Range("C91").Value = "Replica#"
$ws.Range("D91").Value = "Location"

$ws.Range("A93").Value = 1
$ws.Range("B93").Value = 1
$ws.Range("C93").Value = 0

$ws.Range("A94").Value = 2
$ws.Range("B94").Value = 3
$ws.Range("C94").Value = 2

$ws.Range("A95").Value = 3
$ws.Range("B95").Value = 1
$ws.Range("C95").Value = 0

$ws.Range("A96").Value = 4
$ws.Range("B96").Value = 3
$ws.Range("C96").Value = 2

$ws.Range("A97").Value = 5
$ws.Range("B97").Value = 5
$ws.Range("C97").Value = 4

$ws.Range("A98").Value = 6
$ws.Range("B98").Value = 6
$ws.Range("C98").Value = 7

$ws.Range("A99").Value = 7
$ws.Range("B99").Value = 1
$ws.Range("C99").Value = 0

$ws.Range("A100").Value = 8
$ws.Range("B100").Value = 3
$ws.Range("C100").Value = 2
$ws.Range("D100").WrapText = $true

$ws.Range("A101").Value = 9
$ws.Range("B101").Value = 4
$ws.Range("C101").Value = 5

$ws.Range("A102").Value = 10
$ws.Range("B102").Value = 1
$ws.Range("C102").Value = 0

$ws.Range("A103").Value = 11
$ws.Range("B103").Value = 2
$ws.Range("C103").Value = 3

$ws.Range("A104").Value = 12
$ws.Range("B104").Value = 6
$ws.Range("C104").Value = 7

$ws.Range("A105").Value = 13
$ws.Range("B105").Value = 1
$ws.Range("C105").Value = 0

$ws.Range("A106").Value = 14
$ws.Range("B106").Value = 4
$ws.Range("C106").Value = 2

$ws.Range("A107").Value = 15
$ws.Range("B107").Value = 3
$ws.Range("C107").Value = 5
$ws.Range("D107").Value = "*"

$ws.Range("A108").Value = 16
$ws.Range("B108").Value = 6
$ws.Range("C108").Value = 7

$ws.Range("A109").Value = 17
$ws.Range("B109").Value = 1
$ws.Range("C109").Value = 0

$ws.Range("A110").Value = 18
$ws.Range("B110").Value = 3
$ws.Range("C110").Value = 2

$ws.Range("A111").Value = 19
$ws.Range("B111").Value = 4
$ws.Range("C111").Value = 5

$ws.Range("A112").Value = 20
$ws.Range("B112").Value = 6
$ws.Range("C112").Value = 7

$ws.Range("A113").Value = 21
$ws.Range("B113").Value = 3
$ws.Range("C113").Value = 1

$ws.Range("A114").Value = 22
$ws.Range("B114").Value = 2
$ws.Range("C114").Value = 0

$ws.Range("A115").Value = 23
$ws.Range("B115").Value = 4
$ws.Range("C115").Value = 5

$ws.Range("A116").Value = 24
$ws.Range("B116").Value = 6
$ws.Range("C116").Value = 7

$ws.Range("A117").Value = 25
$ws.Range("B117").Value = 3
$ws.Range("C117").Value = 1

$ws.Range("A118").Value = 26
$ws.Range("B118").Value = 2
$ws.Range("C118").Value = 0

$ws.Range("A119").Value = 27
$ws.Range("B119").Value = 4
$ws.Range("C119").Value = 5

$ws.Range("A120").Value = 28
$ws.Range("B120").Value = 6
$ws.Range("C120").Value = 7

$ws.Range("A121").Value = 29
$ws.Range("B121").Value = 3
$ws.Range("C121").Value = 1

$ws.Range("A122").Value = 30
$ws.Range("B122").Value = 2
$ws.Range("C122").Value = 0

$ws.Range("A123").Value = 31
$ws.Range("B123").Value = 4
$ws.Range("C123").Value = 5

$ws.Range("A124").Value = 32
$ws.Range("B124").Value = 6
$ws.Range("C124").Value = 0
$ws.Range("D124").Value = "*"

$ws.Range("A127").Value = "Percentage of local exchanges= 93%"
$ws.Range("D127").Value = "The BigJobs started within 2 minutes of each other"

# --- Page setup: force a portrait pageSetup element to be written ---------
$ws.PageSetup.Orientation = 1

# --- View: scroll down to the new table and select the next empty cell ----
$win = $excel.ActiveWindow
$win.ScrollRow = 102
$win.ScrollColumn = 1
$ws.Range("D128").Select() | Out-Null

# --- Shrink the saved workbook window size to match the authored file -----
$win.Width = 14300
$win.Height = 11800
